$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 3
    "C4" = 1
    "C6" = 0
    "C7" = 2
    "C9" = 1
    "C11" = 0
    "C12" = 8
    "C13" = 0
    "C14" = 3
    "C16" = 0
    "C17" = 0
    "C19" = 3
    "C20" = 2
    "C21" = 0
    "C24" = 0
    "C25" = 6
    "C27" = 0
    "C28" = 1
    "C29" = 1
    "C31" = 4
    "C32" = 0
    "C34" = 0
    "C35" = 0
    "C36" = 0
    "C38" = 0
    "C39" = 0
    "C41" = 0
    "C42" = 1
    "C43" = 0
    "C44" = 1
    "C45" = 0
    "C46" = 0
    "C47" = 0
    "C48" = 1
    "C49" = 3
    "C50" = 0
    "C55" = 0
    "C56" = 6
    "C64" = 1
    "C66" = 6
    "C67" = 1
    "C70" = 2
    "C76" = 3
    "C78" = 6
    "C79" = 2
    "C81" = 0
    "C82" = 4
    "C84" = 4
    "C85" = 10
    "C87" = 4
    "C89" = 2
    "C90" = 4
    "C93" = 5
    "C94" = 0
    "C95" = 1
    "C98" = 9
    "C101" = 1
    "C103" = 3
    "C104" = 4
    "C105" = 15
    "C106" = 8
    "C112" = 0
    "C113" = 1
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
